$wb = $excel.ActiveWorkbook

$wsALC = $wb.Worksheets.Item("ALC")
$wsARM = $wb.Worksheets.Item("ARM")
$wsBSM = $wb.Worksheets.Item("BSM")
$wsCRP = $wb.Worksheets.Item("CRP")
$wsCUL = $wb.Worksheets.Item("CUL")
$wsGSM = $wb.Worksheets.Item("GSM")
$wsLTW = $wb.Worksheets.Item("LTW")
$wsWVR = $wb.Worksheets.Item("WVR")

# --- ALC ---
# row 32
$wsALC.Range("H32").Value = 3402.7273
$wsALC.Range("I32").Value = 1295.6
$wsALC.Range("J32").Value = 5158.6665
$wsALC.Range("K32").Value = 1295.6
$wsALC.Range("L32").Value = 5158.6665
$wsALC.Range("M32").Value = -969.5999999999999
$wsALC.Range("N32").Value = -5810.6665

# row 70
$wsALC.Range("H70").Value = 5549.909
$wsALC.Range("J70").Value = 8591.666999999999
$wsALC.Range("L70").Value = 25775.001
$wsALC.Range("N70").Value = -26315.001

# row 73
$wsALC.Range("H73").Value = 5549.909
$wsALC.Range("J73").Value = 8591.666999999999
$wsALC.Range("L73").Value = 25775.001
$wsALC.Range("N73").Value = -27647.001

# row 98
$wsALC.Range("H98").Value = 592.2258
$wsALC.Range("I98").Value = 592
$wsALC.Range("K98").Value = 592
$wsALC.Range("M98").Value = 906

# row 116
$wsALC.Range("H116").Value = 25800626
$wsALC.Range("I116").Value = 30099614
$wsALC.Range("J116").Value = 6699
$wsALC.Range("K116").Value = 30099614
$wsALC.Range("L116").Value = 6699
$wsALC.Range("M116").Value = -30096172
$wsALC.Range("N116").Value = -13583

# row 122
$wsALC.Range("H122").Value = 592.2258
$wsALC.Range("I122").Value = 592
$wsALC.Range("K122").Value = 1776
$wsALC.Range("M122").Value = 674

# row 132
$wsALC.Range("H132").Value = 14167.021
$wsALC.Range("I132").Value = 2896.389
$wsALC.Range("J132").Value = 21162.586
$wsALC.Range("K132").Value = 8689.167000000001
$wsALC.Range("L132").Value = 63487.758
$wsALC.Range("M132").Value = -6159.167000000001
$wsALC.Range("N132").Value = -68547.758

# row 138
$wsALC.Range("H138").Value = 3027.34
$wsALC.Range("J138").Value = 3452.4558
$wsALC.Range("L138").Value = 10357.3674
$wsALC.Range("N138").Value = -20637.3674

# --- ARM ---
# row 45
$wsARM.Range("H45").Value = 3443.5557
$wsARM.Range("I45").Value = 1999
$wsARM.Range("K45").Value = 1999
$wsARM.Range("M45").Value = -1622

# row 61
$wsARM.Range("H61").Value = 13068.143
$wsARM.Range("I61").Value = 13912.833
$wsARM.Range("K61").Value = 13912.833
$wsARM.Range("M61").Value = -13700.833

# row 74
$wsARM.Range("H74").Value = 1231.6666
$wsARM.Range("I74").Value = 864.8889
$wsARM.Range("K74").Value = 864.8889
$wsARM.Range("M74").Value = 9.111099999999965

# row 77
$wsARM.Range("H77").Value = 1231.6666
$wsARM.Range("I77").Value = 864.8889
$wsARM.Range("K77").Value = 4324.444500000001
$wsARM.Range("M77").Value = 43.55549999999948

# row 122
$wsARM.Range("H122").Value = 5107.68
$wsARM.Range("I122").Value = 2827.7058
$wsARM.Range("K122").Value = 8483.117400000001
$wsARM.Range("M122").Value = -6033.117400000001

# row 132
$wsARM.Range("H132").Value = 15706.435
$wsARM.Range("I132").Value = 27613.773
$wsARM.Range("J132").Value = 4791.375
$wsARM.Range("K132").Value = 82841.319
$wsARM.Range("L132").Value = 14374.125
$wsARM.Range("M132").Value = -80311.319
$wsARM.Range("N132").Value = -19434.125

# row 136
$wsARM.Range("H136").Value = 13068.143
$wsARM.Range("I136").Value = 13912.833
$wsARM.Range("K136").Value = 41738.499
$wsARM.Range("M136").Value = -39188.499

# --- BSM ---
# row 134
$wsBSM.Range("H134").Value = 992
$wsBSM.Range("I134").Value = 942.3333
$wsBSM.Range("K134").Value = 2826.9999
$wsBSM.Range("M134").Value = -291.9998999999998

# --- CRP ---
# row 31
$wsCRP.Range("H31").Value = 5453.0654
$wsCRP.Range("I31").Value = 2179
$wsCRP.Range("J31").Value = 7171.95
$wsCRP.Range("K31").Value = 2179
$wsCRP.Range("L31").Value = 7171.95
$wsCRP.Range("M31").Value = -1884
$wsCRP.Range("N31").Value = -7761.95

# row 34
$wsCRP.Range("H34").Value = 5453.0654
$wsCRP.Range("I34").Value = 2179
$wsCRP.Range("J34").Value = 7171.95
$wsCRP.Range("K34").Value = 2179
$wsCRP.Range("L34").Value = 7171.95
$wsCRP.Range("M34").Value = -1977
$wsCRP.Range("N34").Value = -7575.95

# row 86
$wsCRP.Range("H86").Value = 9692.799999999999
$wsCRP.Range("I86").Value = 9382.777
$wsCRP.Range("K86").Value = 9382.777
$wsCRP.Range("M86").Value = -8259.777

# row 89
$wsCRP.Range("H89").Value = 9692.799999999999
$wsCRP.Range("I89").Value = 9382.777
$wsCRP.Range("K89").Value = 46913.885
$wsCRP.Range("M89").Value = -41297.885

# row 122
$wsCRP.Range("H122").Value = 3362.3215
$wsCRP.Range("I122").Value = 1655.0588
$wsCRP.Range("K122").Value = 4965.1764
$wsCRP.Range("M122").Value = -2515.1764

# row 132
$wsCRP.Range("H132").Value = 9532543
$wsCRP.Range("I132").Value = 11120956
$wsCRP.Range("K132").Value = 33362868
$wsCRP.Range("M132").Value = -33360338

# --- CUL ---
# row 2
$wsCUL.Range("H2").Value = 1613.5834
$wsCUL.Range("I2").Value = 309.7143
$wsCUL.Range("J2").Value = 3439
$wsCUL.Range("K2").Value = 1858.2858
$wsCUL.Range("L2").Value = 20634
$wsCUL.Range("M2").Value = -1745.2858
$wsCUL.Range("N2").Value = -20860

# row 17
$wsCUL.Range("H17").Value = 850
$wsCUL.Range("J17").Value = 850
$wsCUL.Range("L17").Value = 2550
$wsCUL.Range("N17").Value = -2888

# --- GSM ---
# row 132
$wsGSM.Range("H132").Value = 379170.12
$wsGSM.Range("I132").Value = 107829.1
$wsGSM.Range("K132").Value = 323487.3
$wsGSM.Range("M132").Value = -320957.3

# --- LTW ---
# row 16
$wsLTW.Range("H16").Value = 3396.2727
$wsLTW.Range("I16").Value = 1836
$wsLTW.Range("K16").Value = 1836
$wsLTW.Range("M16").Value = -1666

# row 55
$wsLTW.Range("H55").Value = 353.7647
$wsLTW.Range("I55").Value = 203.83333
$wsLTW.Range("J55").Value = 435.54544
$wsLTW.Range("K55").Value = 203.83333
$wsLTW.Range("M55").Value = -30.83332999999999

# row 132
$wsLTW.Range("H132").Value = 3898.6453
$wsLTW.Range("I132").Value = 3809.375
$wsLTW.Range("J132").Value = 4204.7144
$wsLTW.Range("K132").Value = 11428.125
$wsLTW.Range("L132").Value = 12614.1432
$wsLTW.Range("M132").Value = -8898.125
$wsLTW.Range("N132").Value = -17674.1432

# row 136
$wsLTW.Range("H136").Value = 3035.5518
$wsLTW.Range("I136").Value = 2057.6956
$wsLTW.Range("J136").Value = 6784
$wsLTW.Range("K136").Value = 6173.0868
$wsLTW.Range("L136").Value = 20352
$wsLTW.Range("M136").Value = -3623.0868
$wsLTW.Range("N136").Value = -25452

# --- WVR ---
# row 81
$wsWVR.Range("H81").Value = 1751197
$wsWVR.Range("J81").Value = 29322.166
$wsWVR.Range("L81").Value = 58644.332
$wsWVR.Range("N81").Value = -60766.332

# row 84
$wsWVR.Range("H84").Value = 1751197
$wsWVR.Range("J84").Value = 29322.166
$wsWVR.Range("L84").Value = 293221.66
$wsWVR.Range("N84").Value = -303829.66

# row 122
$wsWVR.Range("H122").Value = 3925.5557
$wsWVR.Range("J122").Value = 0
$wsWVR.Range("L122").Value = 0
$wsWVR.Range("N122").ClearContents()

# row 133
$wsWVR.Range("H133").Value = 0
$wsWVR.Range("J133").Value = 0
$wsWVR.Range("L133").Value = 0
$wsWVR.Range("N133").ClearContents()

# row 136
$wsWVR.Range("H136").Value = 7826.551
$wsWVR.Range("I136").Value = 2762.3333
$wsWVR.Range("J136").Value = 10767.064
$wsWVR.Range("K136").Value = 8286.999899999999
$wsWVR.Range("L136").Value = 32301.192
$wsWVR.Range("M136").Value = -5736.999899999999
$wsWVR.Range("N136").Value = -37401.192
